# Appleseed RPA maintenance
# Update the "feb_2021" (column L) figures and roll them into the
# "SFY 2021 Total" (column Q) figures on the crisis_src sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crisis_src")

$updates = @(
    @{Row=3;  L=234; Q=2105}
    @{Row=4;  L=64;  Q=365}
    @{Row=5;  L=298; Q=2470}
    @{Row=7;  L=1;   Q=2}
    @{Row=8;  L=65;  Q=281}
    @{Row=9;  L=8;   Q=24}
    @{Row=10; L=15;  Q=55}
    @{Row=13; L=30;  Q=112}
    @{Row=14; L=119; Q=475}
    @{Row=18; L=22;  Q=132}
    @{Row=19; L=10;  Q=33}
    @{Row=20; L=3;   Q=51}
    @{Row=21; L=67;  Q=336}
    @{Row=22; L=3;   Q=13}
    @{Row=24; L=4;   Q=16}
    @{Row=25; L=1;   Q=1}
    @{Row=26; L=2;   Q=14}
    @{Row=28; L=8;   Q=103}
    @{Row=29; L=12;  Q=84}
)

foreach ($u in $updates) {
    $ws.Range("L$($u.Row)").Value = $u.L
    $ws.Range("Q$($u.Row)").Value = $u.Q
}
